# "fixing some project details"
#
# Three small copy/content fixes on the Projects sheet:
#  - E3  (shortDescription, "My portfolio" row): typo "have" -> "hate"
#  - E17 (shortDescription, "Asema" row): "with my 4 friends" -> "with 4 of my friends"
#  - C18 (title, "Creștem România Împreună" row): fix the un-diacritic'd
#    title "Crestem Romania Impreuna" -> "Creștem România Împreună"
# Plus a cosmetic change of the saved scroll position (sheetView topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "The third redesign of my portfolio... Sure hope I don't hate this one in a few weeks and start from scratch."

$ws.Range("E17").Value = 'Project I cooked up in 48 hours with 4 of my friends at a hackathon organized by the student organization of my university on the topic of "Innovating the current world state through leading-edge & smart technologies".'

$ws.Range("C18").Value = "Creștem România Împreună"

# Best-effort: restore the scroll position the author had when they saved
# (sheetView topLeftCell="B10" -> "C14").
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 3
    $win.ScrollRow = 14
} catch {
    # Scroll/viewport state may not be persisted by this host; ignore.
}
